# Apply COVID-19 country data refresh (12 May 2020, 01:05 -> 01:35)
# Values below are taken directly from the target workbook state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / last-updated timestamp
$ws.Range("A1").Value = 'Datos actualizados a 12 de Mayo de 2020 a las 01:35'

# Estados Unidos - refreshed totals
$ws.Range("B4").Value = 1384526
$ws.Range("C4").Value = 16888
$ws.Range("D4").Value = 260400
$ws.Range("E4").Value = 1042407
$ws.Range("G4").Value = 932
$ws.Range("H4").Value = 81719

# Brasil - refreshed totals
$ws.Range("B11").Value = 168331
$ws.Range("C11").Value = 5632
$ws.Range("D11").Value = 67384
$ws.Range("E11").Value = 89428
$ws.Range("G11").Value = 396
$ws.Range("H11").Value = 11519

# Canada - refreshed totals
$ws.Range("B16").Value = 69981
$ws.Range("C16").Value = 1133
$ws.Range("D16").Value = 32994
$ws.Range("E16").Value = 31994
$ws.Range("G16").Value = 123
$ws.Range("H16").Value = 4993

# Nigeria - refreshed totals
$ws.Range("B64").Value = 4641
$ws.Range("C64").Value = 242
$ws.Range("D64").Value = 902
$ws.Range("E64").Value = 3589
$ws.Range("G64").Value = 7
$ws.Range("H64").Value = 150

# Gabon moves up in the ranking with new totals (was row 118)
$ws.Range("A108").Value = 'Gabon'
$ws.Range("B108").Value = 802
$ws.Range("C108").Value = 141
$ws.Range("D108").Value = 127
$ws.Range("E108").Value = 666
$ws.Range("F108").Value = 1
$ws.Range("H108").Value = 9

# Costa Rica shifts down one place (unchanged totals)
$ws.Range("A109").Value = 'Costa Rica'
$ws.Range("B109").Value = 801
$ws.Range("C109").Value = 9
$ws.Range("D109").Value = 517
$ws.Range("E109").Value = 277
$ws.Range("F109").Value = 6
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 7

# Guinea-Bisau overtakes Burkina Faso, refreshed totals
$ws.Range("A110").Value = 'Guinea-Bisau'
$ws.Range("B110").Value = 761
$ws.Range("C110").Value = 35
$ws.Range("D110").Value = 26
$ws.Range("E110").Value = 732
$ws.Range("H110").Value = 3

# Burkina Faso, refreshed totals
$ws.Range("A111").Value = 'Burkina Faso'
$ws.Range("B111").Value = 760
$ws.Range("C111").Value = 9
$ws.Range("D111").Value = 584
$ws.Range("E111").Value = 126
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 50

# Principado de Andorra shifts down one place (unchanged totals)
$ws.Range("A112").Value = 'Principado de Andorra'
$ws.Range("B112").Value = 755
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 550
$ws.Range("E112").Value = 157
$ws.Range("F112").Value = 14
$ws.Range("H112").Value = 48

# Paraguay shifts down one place (unchanged totals)
$ws.Range("A113").Value = 'Paraguay'
$ws.Range("B113").Value = 724
$ws.Range("C113").Value = 11
$ws.Range("D113").Value = 170
$ws.Range("E113").Value = 544
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 10

# Mali shifts down one place (unchanged totals)
$ws.Range("A114").Value = 'Mali'
$ws.Range("C114").Value = 8
$ws.Range("D114").Value = 377
$ws.Range("E114").Value = 296
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 39

# Crucero shifts down one place (unchanged totals)
$ws.Range("A115").Value = 'Crucero'
$ws.Range("B115").Value = 712
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 651
$ws.Range("E115").Value = 48
$ws.Range("F115").Value = 4
$ws.Range("H115").Value = 13

# Uruguay shifts down one place (unchanged totals)
$ws.Range("A116").Value = 'Uruguay'
$ws.Range("B116").Value = 711
$ws.Range("C116").Value = 4
$ws.Range("D116").Value = 523
$ws.Range("E116").Value = 169
$ws.Range("F116").Value = 8
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 19

# Kenia shifts down one place (unchanged totals)
$ws.Range("A117").Value = 'Kenia'
$ws.Range("B117").Value = 700
$ws.Range("C117").Value = 28
$ws.Range("D117").Value = 251
$ws.Range("E117").Value = 416
$ws.Range("F117").Value = 1
$ws.Range("H117").Value = 33

# Tayikistan shifts down one place (unchanged totals)
$ws.Range("A118").Value = 'Tayikistan'
$ws.Range("C118").Value = 49
$ws.Range("D118").Value = 0
$ws.Range("E118").Value = 640
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 21

# Guyana - refreshed totals
$ws.Range("B162").Value = 109
$ws.Range("C162").Value = 5
$ws.Range("D162").Value = 36
$ws.Range("E162").Value = 63

# Nueva Caledonia overtakes Belice (tie-break swap)
$ws.Range("A192").Value = 'Nueva Caledonia'
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# Belice, unchanged totals
$ws.Range("A193").Value = 'Belice'
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Curazao overtakes Dominica (tie-break swap)
$ws.Range("A198").Value = 'Curazao'
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1

# Dominica, unchanged totals
$ws.Range("A199").Value = 'Dominica'
$ws.Range("D199").Value = 15
$ws.Range("H199").Value = 0
